$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = 1167300
$ws.Range("E12").Value = 1254300
$ws.Range("F12").Value = 1398600
$ws.Range("G12").Value = 1163600
$ws.Range("H12").Value = 704900

$ws.Range("D15").Value = 443800
$ws.Range("E15").Value = 433000
$ws.Range("F15").Value = 432100
$ws.Range("G15").Value = 287900
$ws.Range("H15").Value = 127200

$ws.Range("G21").Value = 1597000
$ws.Range("H21").Value = 1316300

$ws.Range("G83").Value = 508700
$ws.Range("H83").Value = 255100

$ws.Range("D91").Value = -503300
$ws.Range("E91").Value = -463600
$ws.Range("F91").Value = -603800
$ws.Range("G91").Value = -606100
$ws.Range("H91").Value = -565800
$ws.Range("I91").Value = -388400
$ws.Range("J91").Value = -306800
